$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Contest 33: DC vs SRH (row 42) - fill in the points scored against each team
$ws.Range("E42").Value = 100
$ws.Range("H42").Value = 80
$ws.Range("K42").Value = 40
$ws.Range("N42").Value = 20
$ws.Range("Q42").Value = 0
$ws.Range("T42").Value = 60

$excel.CalculateFull()
